$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated loading_percent results for the 380 kV case (rows 2-25)

# Column B
$ws.Cells.Item(2, 2).Value = 10.06655267360777
$ws.Cells.Item(3, 2).Value = 9.754891060528685
$ws.Cells.Item(4, 2).Value = 9.559559341184409
$ws.Cells.Item(5, 2).Value = 9.479081855602898
$ws.Cells.Item(6, 2).Value = 9.465669276293486
$ws.Cells.Item(7, 2).Value = 9.558477384161627
$ws.Cells.Item(8, 2).Value = 9.959985004558233
$ws.Cells.Item(9, 2).Value = 10.71102946378267
$ws.Cells.Item(10, 2).Value = 11.23490512677161
$ws.Cells.Item(11, 2).Value = 11.46612981514405
$ws.Cells.Item(12, 2).Value = 11.55259445849662
$ws.Cells.Item(13, 2).Value = 11.53402257552914
$ws.Cells.Item(14, 2).Value = 11.47326566605273
$ws.Cells.Item(15, 2).Value = 11.43590557922003
$ws.Cells.Item(16, 2).Value = 11.21964473736879
$ws.Cells.Item(17, 2).Value = 11.08510351968687
$ws.Cells.Item(18, 2).Value = 11.00705570513323
$ws.Cells.Item(19, 2).Value = 10.98051858394338
$ws.Cells.Item(20, 2).Value = 11.09949488225769
$ws.Cells.Item(21, 2).Value = 11.49114172752391
$ws.Cells.Item(22, 2).Value = 11.7406928870684
$ws.Cells.Item(23, 2).Value = 11.60811242429921
$ws.Cells.Item(24, 2).Value = 11.09299071637669
$ws.Cells.Item(25, 2).Value = 10.51237233999552

# Column C
$ws.Cells.Item(2, 3).Value = 5.270437362646166
$ws.Cells.Item(3, 3).Value = 5.048455835676474
$ws.Cells.Item(4, 3).Value = 4.907451072836533
$ws.Cells.Item(5, 3).Value = 4.848885965622117
$ws.Cells.Item(6, 3).Value = 4.839097078100396
$ws.Cells.Item(7, 3).Value = 4.906665603393658
$ws.Cells.Item(8, 3).Value = 5.194919490563437
$ws.Cells.Item(9, 3).Value = 5.719852817139502
$ws.Cells.Item(10, 3).Value = 6.07756605640959
$ws.Cells.Item(11, 3).Value = 6.233675111515219
$ws.Cells.Item(12, 3).Value = 6.291800860592578
$ws.Cells.Item(13, 3).Value = 6.279326985122465
$ws.Cells.Item(14, 3).Value = 6.238477176776936
$ws.Cells.Item(15, 3).Value = 6.213325619437475
$ws.Cells.Item(16, 3).Value = 6.067227445526163
$ws.Cells.Item(17, 3).Value = 5.975877426722914
$ws.Cells.Item(18, 3).Value = 5.922714400962665
$ws.Cells.Item(19, 3).Value = 5.904608924609929
$ws.Cells.Item(20, 3).Value = 5.98566634741668
$ws.Cells.Item(21, 3).Value = 6.250502880188065
$ws.Cells.Item(22, 3).Value = 6.417808458406126
$ws.Cells.Item(23, 3).Value = 6.329054291989866
$ws.Cells.Item(24, 3).Value = 5.981242781484742
$ws.Cells.Item(25, 3).Value = 5.582541935349647

# Column D
$ws.Cells.Item(2, 4).Value = 9.193505539852069
$ws.Cells.Item(3, 4).Value = 9.122094076288494
$ws.Cells.Item(4, 4).Value = 9.079538283194617
$ws.Cells.Item(5, 4).Value = 9.062536574309451
$ws.Cells.Item(6, 4).Value = 9.059734436339472
$ws.Cells.Item(7, 4).Value = 9.079307595065428
$ws.Cells.Item(8, 4).Value = 9.168624334730957
$ws.Cells.Item(9, 4).Value = 9.353314830966914
$ws.Cells.Item(10, 4).Value = 9.493850490543336
$ws.Cells.Item(11, 4).Value = 9.558621904290131
$ws.Cells.Item(12, 4).Value = 9.583251850736781
$ws.Cells.Item(13, 4).Value = 9.577943091524199
$ws.Cells.Item(14, 4).Value = 9.560646278498442
$ws.Cells.Item(15, 4).Value = 9.550064268146604
$ws.Cells.Item(16, 4).Value = 9.489632868492672
$ws.Cells.Item(17, 4).Value = 9.452762420166353
$ws.Cells.Item(18, 4).Value = 9.431635791944375
$ws.Cells.Item(19, 4).Value = 9.424497029601071
$ws.Cells.Item(20, 4).Value = 9.456679167766675
$ws.Cells.Item(21, 4).Value = 9.565724140072435
$ws.Cells.Item(22, 4).Value = 9.63757968888312
$ws.Cells.Item(23, 4).Value = 9.599181325311047
$ws.Cells.Item(24, 4).Value = 9.454908184606973
$ws.Cells.Item(25, 4).Value = 9.302428742950745

# Column E
$ws.Cells.Item(2, 5).Value = 13.69638693630328
$ws.Cells.Item(3, 5).Value = 13.62785279992256
$ws.Cells.Item(4, 5).Value = 13.5884402633765
$ws.Cells.Item(5, 5).Value = 13.57306317562707
$ws.Cells.Item(6, 5).Value = 13.57055148619608
$ws.Cells.Item(7, 5).Value = 13.58823009723882
$ws.Cells.Item(8, 5).Value = 13.67221078280898
$ws.Cells.Item(9, 5).Value = 13.85744888539987
$ws.Cells.Item(10, 5).Value = 14.00519182248776
$ws.Cells.Item(11, 5).Value = 14.0747297490087
$ws.Cells.Item(12, 5).Value = 14.10137870078047
$ws.Cells.Item(13, 5).Value = 14.09562559754219
$ws.Cells.Item(14, 5).Value = 14.07691593980308
$ws.Cells.Item(15, 5).Value = 14.06549638830644
$ws.Cells.Item(16, 5).Value = 14.00069265505503
$ws.Cells.Item(17, 5).Value = 13.96152121104678
$ws.Cells.Item(18, 5).Value = 13.93921118410223
$ws.Cells.Item(19, 5).Value = 13.93169579080798
$ws.Cells.Item(20, 5).Value = 13.96566840445366
$ws.Cells.Item(21, 5).Value = 14.08240298119593
$ws.Cells.Item(22, 5).Value = 14.1605313580374
$ws.Cells.Item(23, 5).Value = 14.1186710468187
$ws.Cells.Item(24, 5).Value = 13.96379280231394
$ws.Cells.Item(25, 5).Value = 13.80522453907274

# Column F
$ws.Cells.Item(2, 6).Value = 33.17022212974837
$ws.Cells.Item(3, 6).Value = 33.24478385495912
$ws.Cells.Item(4, 6).Value = 33.29856291471717
$ws.Cells.Item(5, 6).Value = 33.32248484553253
$ws.Cells.Item(6, 6).Value = 33.32657810903505
$ws.Cells.Item(7, 6).Value = 33.29887741635692
$ws.Cells.Item(8, 6).Value = 33.19426843697954
$ws.Cells.Item(9, 6).Value = 33.05277596353501
$ws.Cells.Item(10, 6).Value = 32.98784597866965
$ws.Cells.Item(11, 6).Value = 32.96681784052314
$ws.Cells.Item(12, 6).Value = 32.9600805038408
$ws.Cells.Item(13, 6).Value = 32.96147697770416
$ws.Cells.Item(14, 6).Value = 32.96623898453711
$ws.Cells.Item(15, 6).Value = 32.96931550392777
$ws.Cells.Item(16, 6).Value = 32.98939162259283
$ws.Cells.Item(17, 6).Value = 33.00388861992084
$ws.Cells.Item(18, 6).Value = 33.01302772246032
$ws.Cells.Item(19, 6).Value = 33.016259534318
$ws.Cells.Item(20, 6).Value = 33.00226249205728
$ws.Cells.Item(21, 6).Value = 32.96480699356675
$ws.Cells.Item(22, 6).Value = 32.94747164076765
$ws.Cells.Item(23, 6).Value = 32.95606967884833
$ws.Cells.Item(24, 6).Value = 33.00299515837697
$ws.Cells.Item(25, 6).Value = 33.08421428527517

# Column I
$ws.Cells.Item(2, 9).Value = 22.91223486194552
$ws.Cells.Item(3, 9).Value = 23.01950234982141
$ws.Cells.Item(4, 9).Value = 23.08975856074193
$ws.Cells.Item(5, 9).Value = 23.11949309941544
$ws.Cells.Item(6, 9).Value = 23.12449720394243
$ws.Cells.Item(7, 9).Value = 23.09015509908137
$ws.Cells.Item(8, 9).Value = 22.94830866967536
$ws.Cells.Item(9, 9).Value = 22.70501925660661
$ws.Cells.Item(10, 9).Value = 22.54754697218051
$ws.Cells.Item(11, 9).Value = 22.48053195891911
$ws.Cells.Item(12, 9).Value = 22.45582003356614
$ws.Cells.Item(13, 9).Value = 22.46111258924821
$ws.Cells.Item(14, 9).Value = 22.47848555707318
$ws.Cells.Item(15, 9).Value = 22.48921365759124
$ws.Cells.Item(16, 9).Value = 22.55201959834794
$ws.Cells.Item(17, 9).Value = 22.59173288036356
$ws.Cells.Item(18, 9).Value = 22.61500970308268
$ws.Cells.Item(19, 9).Value = 22.62296548596953
$ws.Cells.Item(20, 9).Value = 22.58746032809999
$ws.Cells.Item(21, 9).Value = 22.47336463462058
$ws.Cells.Item(22, 9).Value = 22.40267427271774
$ws.Cells.Item(23, 9).Value = 22.44004793243636
$ws.Cells.Item(24, 9).Value = 22.58939056247445
$ws.Cells.Item(25, 9).Value = 22.76710104908947

# Column J
$ws.Cells.Item(2, 10).Value = 9.952164006458331
$ws.Cells.Item(3, 10).Value = 9.959068111991911
$ws.Cells.Item(4, 10).Value = 9.964876781477637
$ws.Cells.Item(5, 10).Value = 9.967638875806642
$ws.Cells.Item(6, 10).Value = 9.968121385427537
$ws.Cells.Item(7, 10).Value = 9.964912432268855
$ws.Cells.Item(8, 10).Value = 9.954219023890062
$ws.Cells.Item(9, 10).Value = 9.945684373866975
$ws.Cells.Item(10, 10).Value = 9.946965325269291
$ws.Cells.Item(11, 10).Value = 9.949179269264368
$ws.Cells.Item(12, 10).Value = 9.950251258860835
$ws.Cells.Item(13, 10).Value = 9.950010012993703
$ws.Cells.Item(14, 10).Value = 9.949262786345065
$ws.Cells.Item(15, 10).Value = 9.948835481953559
$ws.Cells.Item(16, 10).Value = 9.94685338476539
$ws.Cells.Item(17, 10).Value = 9.946054591146424
$ws.Cells.Item(18, 10).Value = 9.945748750179058
$ws.Cells.Item(19, 10).Value = 9.945671606346028
$ws.Cells.Item(20, 10).Value = 9.946123733168115
$ws.Cells.Item(21, 10).Value = 9.949475932815268
$ws.Cells.Item(22, 10).Value = 9.953027960627027
$ws.Cells.Item(23, 10).Value = 9.951007975181845
$ws.Cells.Item(24, 10).Value = 9.94609199617747
$ws.Cells.Item(25, 10).Value = 9.946664542727666

# Column K
$ws.Cells.Item(2, 11).Value = 10.24193201029324
$ws.Cells.Item(3, 11).Value = 10.02936680492341
$ws.Cells.Item(4, 11).Value = 9.897778754030858
$ws.Cells.Item(5, 11).Value = 9.843956709547021
$ws.Cells.Item(6, 11).Value = 9.835009701338791
$ws.Cells.Item(7, 11).Value = 9.897053601115925
$ws.Cells.Item(8, 11).Value = 10.16890378237207
$ws.Cells.Item(9, 11).Value = 10.69051175024294
$ws.Cells.Item(10, 11).Value = 11.06294533393328
$ws.Cells.Item(11, 11).Value = 11.22929621004412
$ws.Cells.Item(12, 11).Value = 11.29179090915407
$ws.Cells.Item(13, 11).Value = 11.27835460998384
$ws.Cells.Item(14, 11).Value = 11.23444799646949
$ws.Cells.Item(15, 11).Value = 11.20748734218701
$ws.Cells.Item(16, 11).Value = 11.05200714348507
$ws.Cells.Item(17, 11).Value = 10.95579664484868
$ws.Cells.Item(18, 11).Value = 10.90017303047822
$ws.Cells.Item(19, 11).Value = 10.88129253383465
$ws.Cells.Item(20, 11).Value = 10.96606843026929
$ws.Cells.Item(21, 11).Value = 11.24735840623863
$ws.Cells.Item(22, 11).Value = 11.42826664389221
$ws.Cells.Item(23, 11).Value = 11.33199860085737
$ws.Cells.Item(24, 11).Value = 10.96142552013032
$ws.Cells.Item(25, 11).Value = 10.55103660575688

# Column O
$ws.Cells.Item(2, 15).Value = 25.08334463714854
$ws.Cells.Item(3, 15).Value = 25.17842383260048
$ws.Cells.Item(4, 15).Value = 25.24209190594288
$ws.Cells.Item(5, 15).Value = 25.26936480419913
$ws.Cells.Item(6, 15).Value = 25.27397356552107
$ws.Cells.Item(7, 15).Value = 25.24245434534602
$ws.Cells.Item(8, 15).Value = 25.11502880654694
$ws.Cells.Item(9, 15).Value = 24.90721715727793
$ws.Cells.Item(10, 15).Value = 24.78032356280931
$ws.Cells.Item(11, 15).Value = 24.72822453086912
$ws.Cells.Item(12, 15).Value = 24.70930729869442
$ws.Cells.Item(13, 15).Value = 24.71334532885996
$ws.Cells.Item(14, 15).Value = 24.72665191962974
$ws.Cells.Item(15, 15).Value = 24.7349083533901
$ws.Cells.Item(16, 15).Value = 24.78384178015045
$ws.Cells.Item(17, 15).Value = 24.81530354679371
$ws.Cells.Item(18, 15).Value = 24.83392886955222
$ws.Cells.Item(19, 15).Value = 24.84032591946281
$ws.Cells.Item(20, 15).Value = 24.81189958358447
$ws.Cells.Item(21, 15).Value = 24.72272140595794
$ws.Cells.Item(22, 15).Value = 24.66916975431547
$ws.Cells.Item(23, 15).Value = 24.69731750783354
$ws.Cells.Item(24, 15).Value = 24.81343684114037
$ws.Cells.Item(25, 15).Value = 24.95891737582154
